$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows being appended (USDT daily peg-value history, "new wallet info json format").
# Column A holds date strings (as shared strings, same as all prior rows) and
# column B holds the corresponding USD value for that day.
$dates = @(
    "2024-11-15",
    "2023-08-09",
    "2023-08-05",
    "2023-03-14",
    "2023-03-09",
    "2023-03-07",
    "2023-03-06",
    "2023-03-05",
    "2023-03-04",
    "2023-02-27",
    "2023-02-24",
    "2023-02-23",
    "2023-02-20",
    "2023-02-16",
    "2022-12-13",
    "2022-12-09",
    "2022-12-07",
    "2024-11-16",
    "2024-11-17",
    "2024-11-18",
    "2024-11-20",
    "2024-11-19",
    "2024-11-21",
    "2024-11-22",
    "2024-11-23",
    "2024-11-24"
)

$values = @(
    1,
    1,
    0.9987,
    1.004,
    1,
    0.9999,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    1,
    0.9999,
    1,
    1.001,
    1,
    1.001,
    1.001,
    1.001,
    1.001,
    1.001,
    1.001
)

$startRow = 1021
$endRow = $startRow + $dates.Length - 1

# Write the date strings as formulas that evaluate to plain text first, then
# convert the range to static values via Copy / PasteSpecial(values). Doing it
# this way (rather than Range.Value = "2024-11-15") avoids Excel's usual
# "looks like a date -> store as a date serial number" auto-conversion, so
# these land in the workbook as plain shared-string text, matching the rest
# of column A, with no extra number formatting applied to the cells.
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Formula = '="' + $dates[$i] + '"'
}

$dateRange = $ws.Range("A$startRow" + ":A$endRow")
$dateRange.Copy()
$dateRange.PasteSpecial(-4163)  # xlPasteValues

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
